$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title
$ws.Range("A1").Value = "PRODUCCIÓN DE FRUTAS"

# Month header row (row 2), B2 underlined
$ws.Range("B2").Value = "ENERO"
$ws.Range("B2").Font.Underline = $true
$ws.Range("C2").Value = "FEBRERO"
$ws.Range("D2").Value = "MARZO"
$ws.Range("E2").Value = "ABRIL"
$ws.Range("F2").Value = "MAYO"
$ws.Range("G2").Value = "JUNIO"
$ws.Range("H2").Value = "JULIO"
$ws.Range("I2").Value = "AGOSTO"
$ws.Range("J2").Value = "SEPTIEMBRE"
$ws.Range("K2").Value = "OCTUBRE"
$ws.Range("L2").Value = "NOVIEMBRE"
$ws.Range("M2").Value = "DICIEMBRE"

# Fruit production rows
$ws.Range("A3").Value = "DURAZNO"
$ws.Range("B3").Value = 25871

$ws.Range("A4").Value = "PERA"
$ws.Range("B4").Value = 4589236

$ws.Range("A5").Value = "MANZANA"
$ws.Range("B5").Value = 1458

$ws.Range("A6").Value = "LIMÓN"
$ws.Range("B6").Value = 45879

$ws.Range("A7").Value = "PAPAYA"
$ws.Range("B7").Value = 689521

$ws.Range("A8").Value = "PIÑA"
$ws.Range("B8").Value = 35684

$ws.Range("A9").Value = "NARANJA"
$ws.Range("B9").Value = 59860

$ws.Range("A10").Value = "UVAS"
$ws.Range("B10").Value = 147859

$ws.Range("A11").Value = "ALBARICOQUE"
$ws.Range("B11").Value = 4587

$ws.Range("A12").Value = "SANDÍA"
$ws.Range("B12").Value = 2587793

# Summary labels
$ws.Range("A14").Value = "PROMEDIO X MES"
$ws.Range("A15").Value = "PROMEDIO TOTAL"
$ws.Range("A16").Value = "PRODUCCIÓN X MES"
$ws.Range("A17").Value = "PRODUCCIÓN TOTAL"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 28.22
$ws.Columns.Item(2).ColumnWidth = 11.66

# Selection cursor on B2, matching the saved workbook view
$ws.Range("B2").Select() | Out-Null
